$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.494.39"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.520.21"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'311.03"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'98.85"
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "'0.110"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "'7.21"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "2.907.28"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "'15.32"
$ws.Range("E15").Value = "  -3.62%  "
$ws.Range("D16").Value = "2.543.76"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.805"
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").Value = "42.495.55"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'12.08"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "'69.26"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'240.96"
$ws.Range("E23").Value = "  -2.51%  "
$ws.Range("D24").Value = "'2.84"
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("E25").Value = "  -4.07%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "'38.00"
$ws.Range("E30").Value = "  -6.19%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'5.82"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'156.55"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "'2.71"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'0.0786"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").Value = "'3.15"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'17.48"
$ws.Range("E38").Value = "  -4.71%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'0.117"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").Value = "'4.13"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").Value = "'21.93"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "1.995.47"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").Value = "'8.97"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "2.758.58"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "'78.76"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'71.51"
$ws.Range("E51").Value = "  -2.93%  "
